# Applies the odds updates described by the diff to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 6
$ws.Range("G6").Value = 3.1
$ws.Range("I6").Value = 2.63
$ws.Range("N6").Value = 3.4
$ws.Range("O6").Value = 1.33
$ws.Range("P6").Value = 1.75
$ws.Range("Q6").Value = 2.05
$ws.Range("U6").Value = 13
$ws.Range("W6").Value = 34
$ws.Range("X6").Value = 34
$ws.Range("AE6").Value = 11
$ws.Range("AF6").Value = 12

# Row 16
$ws.Range("G16").Value = 2
$ws.Range("H16").Value = 3.2
$ws.Range("I16").Value = 3.9
$ws.Range("J16").Value = 1.1
$ws.Range("K16").Value = 7
$ws.Range("T16").Value = 6
$ws.Range("U16").Value = 8.5
$ws.Range("Y16").Value = 41

# Row 44
$ws.Range("L44").Value = 1.25
$ws.Range("M44").Value = 3.75
$ws.Range("N44").Value = 1.9
$ws.Range("O44").Value = 1.95
$ws.Range("T44").Value = 11
$ws.Range("U44").Value = 19
$ws.Range("AC44").Value = 51

# Row 45
$ws.Range("G45").Value = 2.87
$ws.Range("H45").Value = 3.35
$ws.Range("I45").Value = 2.25
$ws.Range("M45").Value = 4.9
$ws.Range("O45").Value = 2.35
$ws.Range("R45").Value = 1.39
$ws.Range("S45").Value = 2.57
$ws.Range("T45").Value = 15
$ws.Range("U45").Value = 21
$ws.Range("V45").Value = 10.5
$ws.Range("W45").Value = 40
$ws.Range("X45").Value = 21
$ws.Range("Y45").Value = 20
$ws.Range("Z45").Value = 16
$ws.Range("AA45").Value = 7.3
$ws.Range("AB45").Value = 10
$ws.Range("AE45").Value = 14.5
$ws.Range("AF45").Value = 9
$ws.Range("AG45").Value = 25
$ws.Range("AH45").Value = 16
$ws.Range("AI45").Value = 18.5

# Row 49
$ws.Range("G49").Value = 2.05
$ws.Range("I49").Value = 3.6
$ws.Range("J49").Value = 1.05
$ws.Range("K49").Value = 11
$ws.Range("N49").Value = 1.98
$ws.Range("O49").Value = 1.88
$ws.Range("R49").Value = 1.73
$ws.Range("S49").Value = 2
$ws.Range("T49").Value = 8
$ws.Range("U49").Value = 10
$ws.Range("W49").Value = 19
$ws.Range("Z49").Value = 10
$ws.Range("AB49").Value = 13
$ws.Range("AC49").Value = 41
$ws.Range("AI49").Value = 34
$ws.Range("AJ49").Value = 201
